$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": L6 -> 367.8 (PIEDRA SINTERIZADA for CARAVEDO PAZMIÑO JAHAIRA PAMELA)
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("L6").Value = 367.8
$wsVentasPorGrupo.Range("L20").Value = "1 de 18"

# Sheet "VENTA MENSUAL": F6 (agosto) -> 367.8, F20 (agosto total) -> 367.8
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F6").Value = 367.8
$wsVentaMensual.Range("F20").Value = 367.8

# Sheet "CUMPLIMIENTO MENSUAL": update VENTA/POR CUMPLIR/CUMPLIMIENTO totals
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D2").Value = 68904.48
$wsCumplimiento.Range("E2").Value = -68904.48
$wsCumplimiento.Range("D4").Value = 68904.48
$wsCumplimiento.Range("E4").Value = -53432.9207
$wsCumplimiento.Range("F4").Value = 4.453622202126711
